$d = $word.ActiveDocument

function LastParaInsertPoint() {
    $e = $d.Content.End
    return $d.Range($e - 1, $e - 1)
}

# ---------------------------------------------------------------------
# 1) Blank paragraph right after the "Translation" section's last para
# ---------------------------------------------------------------------
$r = LastParaInsertPoint
$r.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 2) "Build Tools" heading paragraph -- inserted with plain/Normal
#    formatting for now; converted to the Heading1 style at the very
#    end (after everything else exists) so the "next new paragraph
#    inherits the last style touched" behavior of this runtime doesn't
#    leak the Heading1 style (and its stripped contextualSpacing) into
#    the paragraphs that come after it.
# ---------------------------------------------------------------------
$r = LastParaInsertPoint
$r.InsertParagraphAfter()
$r = LastParaInsertPoint
$bmStart = $r.Start
$r.InsertAfter("Build Tools")
$headingParaIndex = $d.Paragraphs.Count

# ---------------------------------------------------------------------
# 3) Blank paragraph
# ---------------------------------------------------------------------
$r = LastParaInsertPoint
$r.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 4) Body paragraph about build tools
# ---------------------------------------------------------------------
$r = LastParaInsertPoint
$r.InsertParagraphAfter()
$r = LastParaInsertPoint
$r.InsertAfter("The two most common build tools for Java are Maven and Gradle.  The most common build tool for Clojure is Leiningen, but Maven can also be used.")

# ---------------------------------------------------------------------
# 5) Blank paragraph
# ---------------------------------------------------------------------
$r = LastParaInsertPoint
$r.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 6) "Choice: " paragraph (bold lead-in + normal text + trailing empty run)
# ---------------------------------------------------------------------
$r = LastParaInsertPoint
$r.InsertParagraphAfter()
$r = LastParaInsertPoint
$boldStart = $r.Start
$r.InsertAfter("Choice: ")
$boldEnd = $d.Paragraphs.Last.Range.End - 1
$r2 = LastParaInsertPoint
$r2.InsertAfter("Leiningen if possible, Maven if not.  Based on our demo that we wrote with Clojure, Leiningen is a very simple tool to use for building projects and automating tests.  It is used for professional Clojure development to deploy to servers as well, so this seems like a good option.  However, if Apache Tika turns out to not play nice with Leiningen, we will use Maven instead.")
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Bold = 1

# Leave a trailing, text-less run (formatting-only) at the end of the
# paragraph -- matches the source document's pattern of a stray empty
# run with just <w:rtl w:val="0"/>. Achieved by splitting off a new
# paragraph mark and then deleting it again, which leaves its run behind.
$r3 = LastParaInsertPoint
$markPos = $r3.Start
$r3.InsertParagraphAfter()
$delRange = $d.Range($markPos, $markPos + 1)
$delRange.Delete()

# ---------------------------------------------------------------------
# Finally: promote the "Build Tools" paragraph to Heading1 and mark its
# bookmark, now that every later paragraph already exists.
# ---------------------------------------------------------------------
$hp = $d.Paragraphs.Item($headingParaIndex)
$hp.Style = "Heading 1"
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("h.itibu7fow21z", $bmRange)

Write-Output "Build Tools section added."
